# Rename the header of column C from "Sign-in Count" to "Sign-in-Count"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Sign-in-Count"
